$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "SCHMITT"
$ws.Range("B7").Value = "Hadrien"
$ws.Range("C7").Value = "21 Bd Maréchal Lyautey, 38000 Grenoble"
$ws.Range("D7").Value = 45.1859198
$ws.Range("E7").Value = 5.731540139877457

$ws.Range("A8").Value = "SCHMITT"
$ws.Range("B8").Value = "Hadrien"
$ws.Range("C8").Value = "21 Bd Maréchal Lyautey, 38000 Grenoble"
$ws.Range("D8").Value = 45.1859198
$ws.Range("E8").Value = 5.731540139877457
